$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 226, shifting rows 226:329 down to 227:330
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with its data
$ws.Cells.Item(226, 1).Value = 7
$ws.Cells.Item(226, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(226, 3).Value = "Ñuble"
$ws.Cells.Item(226, 4).Value = 44825
$ws.Cells.Item(226, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(226, 5).Value = 16
$ws.Cells.Item(226, 6).Value = 100114013
$ws.Cells.Item(226, 7).Value = "Zanahoria"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 120
$ws.Cells.Item(226, 11).Value = 12000
$ws.Cells.Item(226, 12).Value = 13000
$ws.Cells.Item(226, 13).Value = 12500
$ws.Cells.Item(226, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(226, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(226, 16).Value = 625
$ws.Cells.Item(226, 17).Value = 20
$ws.Cells.Item(226, 18).Value = "Hortaliza"
